$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete trailing rows 7, 8 and 9 (delete from the bottom up) ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()

# --- Row 2 ---
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,3).Value = 410
$ws.Cells.Item(2,6).Value = "Matériel"
$ws.Cells.Item(2,9).Value = "Pas Très Importante"
$ws.Cells.Item(2,10).Value = "Freeze de l'écran`n"

# --- Row 3 ---
$ws.Cells.Item(3,1).ClearFormats()
$ws.Cells.Item(3,3).Value = 431
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,6).Value = "Matériel"
$ws.Cells.Item(3,7).Value = "Système (Machine)"
$ws.Cells.Item(3,8).Value = "Code : Python, Java, Html "
$ws.Cells.Item(3,9).Value = "Pas Très Importante"
$ws.Cells.Item(3,10).Value = "Freeze pendant alt tab"

# --- Row 4 ---
$ws.Cells.Item(4,1).ClearFormats()
$ws.Cells.Item(4,3).Value = 2095
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,6).Value = "Perturbation"
$ws.Cells.Item(4,8).Value = "Code : Python, Java, Html "
$ws.Cells.Item(4,9).Value = "Pas Très Importante"
$ws.Cells.Item(4,10).Value = "Quelqu'un à crier dans la rue`n"

# --- Row 5 ---
$ws.Cells.Item(5,1).ClearFormats()
$ws.Cells.Item(5,3).Value = 2136
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,6).ClearContents()
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(5,8).ClearContents()
$ws.Cells.Item(5,9).ClearContents()
$ws.Cells.Item(5,10).ClearContents()
$ws.Cells.Item(5,11).ClearContents()
$ws.Cells.Item(5,14).ClearContents()

# --- Row 6 ---
$ws.Cells.Item(6,1).ClearFormats()
$ws.Cells.Item(6,3).Value = 2138
$ws.Cells.Item(6,4).Value = 4
$ws.Cells.Item(6,5).Value = 4
